$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-5 with the new TPM-derived values and re-point the
# "Target cluster" (column D) to Resolving-Mac, matching the refreshed
# NATMI output for Icam2-Itgam (YoungD7).

$ws.Range("A2").Value = "ECs"
$ws.Range("D2").Value = "Resolving-Mac"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 52.08999
$ws.Range("H2").Value = 156.26997
$ws.Range("I2").Value = 0.9401105828221099
$ws.Range("J2").Value = 0.9401105828221098
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 35.68243999999999
$ws.Range("N2").Value = 107.04732
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 1858.6979427756
$ws.Range("R2").Value = 16728.2814849804
$ws.Range("S2").Value = 0.9401105828221099
$ws.Range("T2").Value = 0.9401105828221098

$ws.Range("A3").Value = "FAPs"
$ws.Range("D3").Value = "Resolving-Mac"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.545858
$ws.Range("H3").Value = 1.637574
$ws.Range("I3").Value = 0.009851545038079508
$ws.Range("J3").Value = 0.009851545038079508
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 35.68243999999999
$ws.Range("N3").Value = 107.04732
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 19.47754533351999
$ws.Range("R3").Value = 175.29790800168
$ws.Range("S3").Value = 0.009851545038079508
$ws.Range("T3").Value = 0.009851545038079508

$ws.Range("A4").Value = "MuSCs"
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 1.744358333333333
$ws.Range("H4").Value = 5.233075
$ws.Range("I4").Value = 0.03148185917103467
$ws.Range("J4").Value = 0.03148185917103467
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 35.68243999999999
$ws.Range("N4").Value = 107.04732
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 62.24296156766666
$ws.Range("R4").Value = 560.186654109
$ws.Range("S4").Value = 0.03148185917103467
$ws.Range("T4").Value = 0.03148185917103467

$ws.Range("A5").Value = "Resolving-Mac"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.028158333333334
$ws.Range("H5").Value = 3.084475
$ws.Range("I5").Value = 0.01855601296877595
$ws.Range("J5").Value = 0.01855601296877594
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 35.68243999999999
$ws.Range("N5").Value = 107.04732
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 36.68719803966667
$ws.Range("R5").Value = 330.184782357
$ws.Range("S5").Value = 0.01855601296877595
$ws.Range("T5").Value = 0.01855601296877594

# Remove the old rows 6-9 (MuSCs/Resolving-Mac x ECs/Resolving-Mac pairs)
# now that the data set only keeps the Resolving-Mac target rows.
$ws.Rows("6:9").Delete() | Out-Null
